$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2451187.8
$ws.Range("I9").Value = 3676615.2
$ws.Range("K9").Value = 3676615.2
$ws.Range("M9").Value = -3676446.2
$ws.Range("H12").Value = 902.9231
$ws.Range("I12").Value = 521.63635
$ws.Range("K12").Value = 521.63635
$ws.Range("M12").Value = -351.63635
$ws.Range("H70").Value = 6443.6665
$ws.Range("I70").Value = 6039
$ws.Range("J70").Value = 6949.5
$ws.Range("K70").Value = 18117
$ws.Range("L70").Value = 20848.5
$ws.Range("M70").Value = -17847
$ws.Range("N70").Value = -21388.5
$ws.Range("H73").Value = 6443.6665
$ws.Range("I73").Value = 6039
$ws.Range("J73").Value = 6949.5
$ws.Range("K73").Value = 18117
$ws.Range("L73").Value = 20848.5
$ws.Range("M73").Value = -17181
$ws.Range("N73").Value = -22720.5
$ws.Range("H113").Value = 3910.2727
$ws.Range("I113").Value = 3082.6667
$ws.Range("K113").Value = 3082.6667
$ws.Range("M113").Value = 171.3332999999998
$ws.Range("H132").Value = 6121.24
$ws.Range("I132").Value = 6272.125
$ws.Range("K132").Value = 18816.375
$ws.Range("M132").Value = -16286.375
$ws.Range("H137").Value = 2944651.5
$ws.Range("I137").Value = 4168625.8
$ws.Range("J137").Value = 7113.4
$ws.Range("K137").Value = 12505877.4
$ws.Range("L137").Value = 21340.2
$ws.Range("M137").Value = -12503327.4
$ws.Range("N137").Value = -26440.2
$ws.Range("H141").Value = 2251.7273
$ws.Range("I141").Value = 2134.0527
$ws.Range("J141").Value = 2997
$ws.Range("K141").Value = 6402.158100000001
$ws.Range("L141").Value = 8991
$ws.Range("M141").Value = -1222.158100000001
$ws.Range("N141").Value = -19351

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1931.909
$ws.Range("J6").Value = 2111
$ws.Range("L6").Value = 2111
$ws.Range("N6").Value = -2457
$ws.Range("H63").Value = 2103.5
$ws.Range("I63").Value = 2103.5
$ws.Range("K63").Value = 2103.5
$ws.Range("M63").Value = -1417.5
$ws.Range("H66").Value = 2103.5
$ws.Range("I66").Value = 2103.5
$ws.Range("K66").Value = 10517.5
$ws.Range("M66").Value = -7085.5
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17247772
$ws.Range("I20").Value = 20840396
$ws.Range("K20").Value = 20840396
$ws.Range("M20").Value = -20840149
$ws.Range("H86").Value = 2848.5
$ws.Range("I86").Value = 2202.5386
$ws.Range("J86").Value = 3611.9092
$ws.Range("K86").Value = 2202.5386
$ws.Range("L86").Value = 3611.9092
$ws.Range("M86").Value = -1079.5386
$ws.Range("N86").Value = -5857.9092
$ws.Range("H89").Value = 2848.5
$ws.Range("I89").Value = 2202.5386
$ws.Range("J89").Value = 3611.9092
$ws.Range("K89").Value = 11012.693
$ws.Range("L89").Value = 18059.546
$ws.Range("M89").Value = -5396.692999999999
$ws.Range("N89").Value = -29291.546
$ws.Range("H99").Value = 4077.4
$ws.Range("I99").Value = 3221.75
$ws.Range("K99").Value = 3221.75
$ws.Range("M99").Value = -1723.75
$ws.Range("H135").Value = 69449.75
$ws.Range("J135").Value = 69449.75
$ws.Range("L135").Value = 69449.75
$ws.Range("N135").Value = -79589.75
$ws.Range("H137").Value = 64999
$ws.Range("J137").Value = 64999
$ws.Range("L137").Value = 64999
$ws.Range("N137").Value = -75199

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2019714.2
$ws.Range("I31").Value = 3023.0244
$ws.Range("J31").Value = 5957064
$ws.Range("K31").Value = 3023.0244
$ws.Range("L31").Value = 5957064
$ws.Range("M31").Value = -2728.0244
$ws.Range("N31").Value = -5957654
$ws.Range("H34").Value = 2019714.2
$ws.Range("I34").Value = 3023.0244
$ws.Range("J34").Value = 5957064
$ws.Range("K34").Value = 3023.0244
$ws.Range("L34").Value = 5957064
$ws.Range("M34").Value = -2821.0244
$ws.Range("N34").Value = -5957468
$ws.Range("H86").Value = 7119.778
$ws.Range("I86").Value = 7163.1333
$ws.Range("K86").Value = 7163.1333
$ws.Range("M86").Value = -6040.1333
$ws.Range("H89").Value = 7119.778
$ws.Range("I89").Value = 7163.1333
$ws.Range("K89").Value = 35815.6665
$ws.Range("M89").Value = -30199.6665
$ws.Range("H99").Value = 3599.6
$ws.Range("J99").Value = 4666.6665
$ws.Range("L99").Value = 4666.6665
$ws.Range("N99").Value = -7662.6665
$ws.Range("H126").Value = 3599.6
$ws.Range("J126").Value = 4666.6665
$ws.Range("L126").Value = 13999.9995
$ws.Range("N126").Value = -18939.9995
$ws.Range("H134").Value = 2579.8718
$ws.Range("I134").Value = 2355.0908
$ws.Range("K134").Value = 7065.2724
$ws.Range("M134").Value = -4530.2724

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2167423.5
$ws.Range("J4").Value = 14070202
$ws.Range("L4").Value = 42210606
$ws.Range("N4").Value = -42210830
$ws.Range("H22").Value = 3621.2144
$ws.Range("J22").Value = 4066.4167
$ws.Range("L22").Value = 12199.2501
$ws.Range("N22").Value = -12537.2501
$ws.Range("H27").Value = 3621.2144
$ws.Range("J27").Value = 4066.4167
$ws.Range("L27").Value = 12199.2501
$ws.Range("N27").Value = -12403.2501
$ws.Range("H32").Value = 116669980
$ws.Range("J32").Value = 18522202
$ws.Range("L32").Value = 55566606
$ws.Range("N32").Value = -55567172
$ws.Range("H55").Value = 4843.36
$ws.Range("I55").Value = 2342
$ws.Range("J55").Value = 5816.1113
$ws.Range("K55").Value = 7026
$ws.Range("L55").Value = 17448.3339
$ws.Range("M55").Value = -6849
$ws.Range("N55").Value = -17802.3339
$ws.Range("H59").Value = 831.6667
$ws.Range("I59").Value = 747.5
$ws.Range("J59").Value = 1000
$ws.Range("K59").Value = 2242.5
$ws.Range("L59").Value = 3000
$ws.Range("M59").Value = -1702.5
$ws.Range("N59").Value = -4080
$ws.Range("H107").Value = 4993.4443
$ws.Range("J107").Value = 4665.5
$ws.Range("L107").Value = 13996.5
$ws.Range("N107").Value = -17836.5
$ws.Range("H122").Value = 1581.9048
$ws.Range("I122").Value = 667.3333
$ws.Range("J122").Value = 1734.3334
$ws.Range("K122").Value = 6005.9997
$ws.Range("L122").Value = 15609.0006
$ws.Range("M122").Value = -3555.9997
$ws.Range("N122").Value = -20509.0006
$ws.Range("H131").Value = 5821371.5
$ws.Range("J131").Value = 2037.75
$ws.Range("L131").Value = 6113.25
$ws.Range("N131").Value = -16193.25
$ws.Range("H134").Value = 2675.0667
$ws.Range("I134").Value = 1125.1111
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 3375.3333
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 1694.6667
$ws.Range("N134").Value = -25140

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 425
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 500000
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H122").Value = 11024.782
$ws.Range("I122").Value = 9417.556
$ws.Range("K122").Value = 28252.668
$ws.Range("M122").Value = -25802.668
$ws.Range("H141").Value = 69249
$ws.Range("J141").Value = 69249
$ws.Range("L141").Value = 69249
$ws.Range("N141").Value = -79609

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 30000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 30000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 30000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -30580
$ws.Range("H41").Value = 18891.6
$ws.Range("J41").Value = 18864.5
$ws.Range("L41").Value = 18864.5
$ws.Range("N41").Value = -19644.5
